$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.369.12'
$ws.Range('D3').Value = '3.483.62'
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.52'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.31'
$ws.Range('E6').Value = '  +2.44%  '
$ws.Range('D7').Value = '3.481.72'
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.486'
$ws.Range('E9').Value = '  -0.67%  '
$ws.Range('E10').Value = '  +0.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.19'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.375'
$ws.Range('E12').Value = '  -2.60%  '
$ws.Range('D13').Value = '4.086.95'
$ws.Range('E13').Value = '  +0.29%  '
$ws.Range('E14').Value = '  +1.31%  '
$ws.Range('E15').Value = '  +1.26%  '
$ws.Range('D16').Value = '3.491.10'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '64.428.35'
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.21'
$ws.Range('E18').Value = '  -8.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.99'
$ws.Range('E19').Value = '  +1.60%  '
$ws.Range('E20').Value = '  -0.07%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.71'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '386.62'
$ws.Range('E22').Value = '  -1.36%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.565'
$ws.Range('E23').Value = '  -1.35%  '
$ws.Range('D24').Value = '3.626.21'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '73.98'
$ws.Range('E25').Value = '  +1.40%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('E27').Value = '  +4.79%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.02%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.44'
$ws.Range('E29').Value = '  +1.32%  '
$ws.Range('E30').Value = '  -0.98%  '
$ws.Range('E31').Value = '  -0.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.19'
$ws.Range('E32').Value = '  +0.87%  '
$ws.Range('D33').Value = '3.505.32'
$ws.Range('E33').Value = '  +0.93%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.146'
$ws.Range('E35').Value = '  +1.82%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.37'
$ws.Range('E36').Value = '  -1.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.28'
$ws.Range('E37').Value = '  +2.15%  '
$ws.Range('B38').Value = 'Aptos'
$ws.Range('C38').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.83'
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('B39').Value = 'ImmutableX'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.54'
$ws.Range('E39').Value = '  -1.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '162.86'
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0777'
$ws.Range('E41').Value = '  -2.93%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.802'
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '25.67'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.83'
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('E46').Value = '  +1.77%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.40'
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.66'
$ws.Range('E48').Value = '  +2.32%  '
$ws.Range('D49').Value = '2.469.60'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.902'
$ws.Range('E51').Value = '  +1.99%  '
